$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row right before the existing row 375 (old date 2021-02-08 / 44236),
# shifting the existing rows 375-487 down to 376-488.
$ws.Rows("375:375").Insert()

$ws.Cells.Item(375, 1).Value = 4
$ws.Cells.Item(375, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(375, 3).Value = "Los Lagos"
$ws.Cells.Item(375, 4).Value = 45120
$ws.Cells.Item(375, 5).Value = 10
$ws.Cells.Item(375, 6).Value = 100112003
$ws.Cells.Item(375, 7).Value = "Ajo"
$ws.Cells.Item(375, 8).Value = "Chino"
$ws.Cells.Item(375, 9).Value = "Primera"
$ws.Cells.Item(375, 10).Value = 80
$ws.Cells.Item(375, 11).Value = 22000
$ws.Cells.Item(375, 12).Value = 22000
$ws.Cells.Item(375, 13).Value = 22000
$ws.Cells.Item(375, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(375, 15).Value = "China"
$ws.Cells.Item(375, 16).Value = 2200
$ws.Cells.Item(375, 17).Value = 10
$ws.Cells.Item(375, 18).Value = "Hortaliza"

# Insert a second new data row before what is now row 485 (old date 2022-03-13 / 44662),
# shifting the current rows 485-488 down to 486-489.
$ws.Rows("485:485").Insert()

$ws.Cells.Item(485, 1).Value = 4
$ws.Cells.Item(485, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(485, 3).Value = "Los Lagos"
$ws.Cells.Item(485, 4).Value = 45121
$ws.Cells.Item(485, 5).Value = 10
$ws.Cells.Item(485, 6).Value = 100112003
$ws.Cells.Item(485, 7).Value = "Ajo"
$ws.Cells.Item(485, 8).Value = "Chino"
$ws.Cells.Item(485, 9).Value = "Primera"
$ws.Cells.Item(485, 10).Value = 240
$ws.Cells.Item(485, 11).Value = 22000
$ws.Cells.Item(485, 12).Value = 22000
$ws.Cells.Item(485, 13).Value = 22000
$ws.Cells.Item(485, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(485, 15).Value = "China"
$ws.Cells.Item(485, 16).Value = 2200
$ws.Cells.Item(485, 17).Value = 10
$ws.Cells.Item(485, 18).Value = "Hortaliza"
